$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are a data engineer at a health tech startup. The company wants to move its on-premises healthcare records to GCP. You are tasked with designing a managed, cloud-based data processing solution that can ingest, prepare, transform, and move this data to GCP storage solutions at scale.Which GCP service should you use?",
        "ques_type": 2,
        "options": [
            "Dataflow",
            "Cloud Storage",
            "BigQuery",
            "Cloud Firestore"
        ],
        "score": "Dataflow"
    },
    {
        "title": "You work for a manufacturing firm that has raw supply chain data which requires cleansing, transformation, and enrichment before analytics. The firm needs a GCP service that integrates seamlessly with BigQuery and provides an intuitive interface for transformations. Your task is to devise the required setup.Which GCP service should you use?",
        "ques_type": 2,
        "options": [
            "Dataprep",
            "Cloud Pub/Sub",
            "Cloud Functions",
            "Dataproc"
        ],
        "score": "Dataprep"
    },
    {
        "title": "You are a data analyst at a fintech startup that wants to gain insights from its growing transaction data. The company is already using GCP for its infrastructure, and you have been tasked with creating dashboards that visually represent financial trends.What should you do?",
        "ques_type": 2,
        "options": [
            "Utilize Looker Studio to create interactive dashboards and reports.",
            "Run SQL queries on BigQuery and plot charts.",
            "Implement Cloud Endpoints to fetch and display data.",
            "Analyze the raw data using Dataproc."
        ],
        "score": "Utilize Looker Studio to create interactive dashboards and reports."
    },
    {
        "title": "You are the lead developer for a social media app with a global user base. The app's backend is experiencing high latency issues during peak times. You are tasked with choosing a GCP service that can handle millions of reads/writes per second with low latency.Which GCP service should you use?",
        "ques_type": 2,
        "options": [
            "Cloud Spanner",
            "Cloud SQL",
            "Cloud Bigtable",
            "Cloud Storage"
        ],
        "score": "Cloud Spanner"
    }
]
'@

# Remove trailing newline added by here-string
$newText = $newText.TrimEnd("`r", "`n")

# Clear the old A1 (value 0, bold/bordered style) and A2 contents
$ws.Cells.Item(1, 1).ClearContents()
$ws.Cells.Item(1, 1).ClearFormats()
$ws.Cells.Item(2, 1).ClearContents()
$ws.Cells.Item(2, 1).ClearFormats()

# Set new content in A1 with default style
$ws.Range("A1").Value = $newText
